$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.573811
$ws.Range("H2").Value = 13.721433
$ws.Range("I2").Value = 0.1659009079913533
$ws.Range("J2").Value = 0.1659009079913533
$ws.Range("M2").Value = 0.8317113333333332
$ws.Range("N2").Value = 2.495134
$ws.Range("O2").Value = 0.0263454906755698
$ws.Range("P2").Value = 0.0263454906755698
$ws.Range("Q2").Value = 3.804090445224666
$ws.Range("R2").Value = 34.236814007022
$ws.Range("S2").Value = 0.004370740824554762
$ws.Range("T2").Value = 0.004370740824554763
$ws.Range("G3").Value = 4.573811
$ws.Range("H3").Value = 13.721433
$ws.Range("I3").Value = 0.1659009079913533
$ws.Range("J3").Value = 0.1659009079913533
$ws.Range("O3").Value = 0.6529848313028861
$ws.Range("P3").Value = 0.6529848313028862
$ws.Range("Q3").Value = 94.28609200053268
$ws.Range("R3").Value = 848.5748280047941
$ws.Range("S3").Value = 0.1083307764177295
$ws.Range("T3").Value = 0.1083307764177295
$ws.Range("G4").Value = 4.573811
$ws.Range("H4").Value = 13.721433
$ws.Range("I4").Value = 0.1659009079913533
$ws.Range("J4").Value = 0.1659009079913533
$ws.Range("M4").Value = 10.12334933333333
$ws.Range("N4").Value = 30.370048
$ws.Range("O4").Value = 0.3206696780215441
$ws.Range("P4").Value = 0.3206696780215441
$ws.Range("Q4").Value = 46.30228653764267
$ws.Range("R4").Value = 416.720578838784
$ws.Range("S4").Value = 0.05319939074906908
$ws.Range("T4").Value = 0.05319939074906908
$ws.Range("I5").Value = 0.5322852674812913
$ws.Range("J5").Value = 0.5322852674812913
$ws.Range("M5").Value = 0.8317113333333332
$ws.Range("N5").Value = 2.495134
$ws.Range("O5").Value = 0.0263454906755698
$ws.Range("P5").Value = 0.0263454906755698
$ws.Range("Q5").Value = 12.20524543642022
$ws.Range("R5").Value = 109.847208927782
$ws.Range("S5").Value = 0.01402331655117154
$ws.Range("T5").Value = 0.01402331655117154
$ws.Range("I6").Value = 0.5322852674812913
$ws.Range("J6").Value = 0.5322852674812913
$ws.Range("O6").Value = 0.6529848313028861
$ws.Range("P6").Value = 0.6529848313028862
$ws.Range("S6").Value = 0.3475742055912826
$ws.Range("T6").Value = 0.3475742055912827
$ws.Range("I7").Value = 0.5322852674812913
$ws.Range("J7").Value = 0.5322852674812913
$ws.Range("M7").Value = 10.12334933333333
$ws.Range("N7").Value = 30.370048
$ws.Range("O7").Value = 0.3206696780215441
$ws.Range("P7").Value = 0.3206696780215441
$ws.Range("Q7").Value = 148.5587105766116
$ws.Range("R7").Value = 1337.028395189504
$ws.Range("S7").Value = 0.1706877453388372
$ws.Range("T7").Value = 0.1706877453388372
$ws.Range("G8").Value = 8.320867
$ws.Range("H8").Value = 24.962601
$ws.Range("I8").Value = 0.3018138245273554
$ws.Range("J8").Value = 0.3018138245273554
$ws.Range("M8").Value = 0.8317113333333332
$ws.Range("N8").Value = 2.495134
$ws.Range("O8").Value = 0.0263454906755698
$ws.Range("P8").Value = 0.0263454906755698
$ws.Range("Q8").Value = 6.920559387059332
$ws.Range("R8").Value = 62.28503448353399
$ws.Range("S8").Value = 0.0079514332998435
$ws.Range("T8").Value = 0.007951433299843503
$ws.Range("G9").Value = 8.320867
$ws.Range("H9").Value = 24.962601
$ws.Range("I9").Value = 0.3018138245273554
$ws.Range("J9").Value = 0.3018138245273554
$ws.Range("O9").Value = 0.6529848313028861
$ws.Range("P9").Value = 0.6529848313028862
$ws.Range("Q9").Value = 171.5291758855353
$ws.Range("R9").Value = 1543.762582969818
$ws.Range("S9").Value = 0.197079849293874
$ws.Range("T9").Value = 0.1970798492938741
$ws.Range("G10").Value = 8.320867
$ws.Range("H10").Value = 24.962601
$ws.Range("I10").Value = 0.3018138245273554
$ws.Range("J10").Value = 0.3018138245273554
$ws.Range("M10").Value = 10.12334933333333
$ws.Range("N10").Value = 30.370048
$ws.Range("O10").Value = 0.3206696780215441
$ws.Range("P10").Value = 0.3206696780215441
$ws.Range("Q10").Value = 84.23504339720533
$ws.Range("R10").Value = 758.115390574848
$ws.Range("S10").Value = 0.09678254193363785
$ws.Range("T10").Value = 0.09678254193363786
